# manual data load test.xlsx - update callback data
# Removed redundant calls when creating update callbacks:
#  - clear the stale IMPORT STATUS result cell (K2) from a prior run
#  - change the "New Sub1" name-add row (row 4) to add a synonym instead of
#    a duplicate/placeholder entry
#  - add a duplicate-name batch row (row 6) and a new Nicotine name row (row 8)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$wsCreate = $wb.Worksheets.Item("Substance Creation")

# Clear the leftover IMPORT STATUS value from row 2 (K2) - the callback
# no longer redundantly stamps a result before the batch is resubmitted.
$ws.Range("K2").ClearContents()

# Update row 4 ("New Sub1" / cn name) to reference the new synonym data.
$ws.Range("D4").Value = "Another synonym"
$ws.Range("I4").Value = "blah again"
$ws.Range("J4").Value = "Adding another name"

# New row 6: duplicate Levothyroxine synonym-add batch entry (new UUID).
$ws.Range("B6").Value = "5259a1b9-1627-4431-8fce-5bceb15ce47c"
$ws.Range("C6").Value = "Levothyroxine"
$ws.Range("D6").Value = "3,5,3',5'-Tetraiodo-L-thyronine"
$ws.Range("E6").Value = "sys"
$ws.Range("F6").Value = "Y"
$ws.Range("G6").Value = "WEBSITE"
$ws.Range("H6").Value = "Levo"
$ws.Range("I6").Value = "https://en.wikipedia.org/wiki/Levothyroxine"
$ws.Range("J6").Value = "Adding a name"
$ws.Range("B6:J6").Style = $ws.Range("C2").Style

# New row 8: Nicotine name-add batch entry.
$ws.Range("B8").Value = "a5fe114a-d4e9-4f13-bc75-1503423a7a7d"
$ws.Range("C8").Value = "Nicotine "
$ws.Range("D8").Value = "(-)-3-(1-Methyl-2-pyrrolidyl)pyridine"
$ws.Range("F8").Value = "Y"
$ws.Range("G8").Value = "WEBSITE"
$ws.Range("H8").Value = "nicotine"
$ws.Range("I8").Value = "https://chem.nlm.nih.gov/chemidplus/name/nicotine"
$ws.Range("J8").Value = "adding a name"
$ws.Range("E8").Value = "SYS"
$ws.Range("B8:J8").Style = $ws.Range("C2").Style

# Refresh the view state: select row 3 on "Substance Creation" (no longer
# scrolled to column C), then return focus to "Sheet1" with A8 selected.
$wsCreate.Rows.Item(3).Select() | Out-Null
$ws.Select() | Out-Null
$ws.Range("A8").Select() | Out-Null
